# kröfulisti með c kröfu
#
# The "notendur" (user group) column on several requirement rows is
# normalised to the single value "Allir" (everyone), and a brand-new
# priority-C requirement row is filled in at the bottom of the "kröfur"
# sheet (it previously existed only as an empty placeholder row).
#
# Shared strings are appended/removed automatically by the engine as a
# side effect of the cell-value edits below (matching the order the
# author must have made them in: the new row first, then the "Allir"
# normalisation), so no direct shared-string bookkeeping is needed here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kröfur")

# --- New row 38 / sheet row 41: "C" priority requirement about seeing
#     how many pilots are certified per plane type. ---
$ws.Range("C41").Value = "Það þarf að vera hægt að sjá lista yfir allar flugvélategund og hve margir flugmenn hafa réttindi á viðkomandi tegund"
$ws.Range("E41").Value = "C"
$ws.Range("D41").Value = "Mannari"

# --- Normalise the "user group(s)" column to "Allir" on the rows that
#     used to read "Skipuleggjari", "Skipuleggjari og Mannari" or the
#     lowercase "allir". ---
$ws.Range("D15").Value = "Allir"
$ws.Range("D16").Value = "Allir"
$ws.Range("D33").Value = "Allir"
$ws.Range("D34").Value = "Allir"
$ws.Range("D35").Value = "Allir"
$ws.Range("D37").Value = "Allir"
$ws.Range("D38").Value = "Allir"
$ws.Range("D39").Value = "Allir"
$ws.Range("D40").Value = "Allir"

# --- View state: zoom out to 90%, scroll back to the top-left and move
#     the selection to I3. ---
$ws.Activate()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("I3").Select() | Out-Null
$win.Zoom = 90
